$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.025.92'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.863.41'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.24'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5096'
$ws.Range('E7').Value = '  +1.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3841'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08317'
$ws.Range('E9').Value = '  -7.02%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.52'
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.228'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.61'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').Value = '1.856.32'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.221'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001099'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.91'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06632'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.71'
$ws.Range('E20').Value = '  -2.56%  '
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.044'
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').Value = '28.053.13'
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.11'
$ws.Range('E24').Value = '  -3.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.233'
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.549'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('D27').Value = '2.074.61'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '158.00'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.80'
$ws.Range('E30').Value = '  -1.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1055'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.927'
$ws.Range('E33').Value = '  +5.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.591'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.454'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02421'
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06533'
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2174'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.206'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6477'
$ws.Range('E40').Value = '  +1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.226'
$ws.Range('E41').Value = '  -4.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.968'
$ws.Range('E42').Value = '  +1.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.22'
$ws.Range('E43').Value = '  -2.32%  '
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.10'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.667'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('E48').Value = '  +1.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.208'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '120.16'
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.30'
